$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.116.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.52"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5245"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06343"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.34"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07804"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.500"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.657.60"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5469"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8195"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.35"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.109.54"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.579"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.21"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.06"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.028"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "141.99"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1238"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.235"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.14"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.432"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05897"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.516"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.248"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.584"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9505"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.785"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5672"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01616"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.820"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8474"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.028.84"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.67"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.800.12"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4304"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05167"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.862"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.469"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09699"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.25%  "
